$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 155
$ws.Cells.Item(3, 6).Value = 929
$ws.Cells.Item(4, 6).Value = 1090
$ws.Cells.Item(5, 6).Value = 1554
$ws.Cells.Item(6, 6).Value = 339
$ws.Cells.Item(7, 6).Value = 697
$ws.Cells.Item(8, 6).Value = 12532
$ws.Cells.Item(9, 6).Value = 2210
$ws.Cells.Item(11, 6).Value = 275
$ws.Cells.Item(12, 6).Value = 16532
$ws.Cells.Item(14, 6).Value = 1252
$ws.Cells.Item(15, 6).Value = 238
$ws.Cells.Item(17, 6).Value = 801
$ws.Cells.Item(19, 6).Value = 319
$ws.Cells.Item(21, 6).Value = 788
$ws.Cells.Item(22, 6).Value = 4475
$ws.Cells.Item(23, 6).Value = 1158
$ws.Cells.Item(24, 6).Value = 885
$ws.Cells.Item(25, 6).Value = 13
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(29, 6).Value = 1104
$ws.Cells.Item(30, 6).Value = 59
$ws.Cells.Item(31, 6).Value = 120
$ws.Cells.Item(32, 6).Value = 283
$ws.Cells.Item(36, 6).Value = 25
$ws.Cells.Item(37, 6).Value = 4502
$ws.Cells.Item(39, 6).Value = 4614
$ws.Cells.Item(40, 6).Value = 5581
$ws.Cells.Item(43, 6).Value = 92
$ws.Cells.Item(44, 6).Value = 179
$ws.Cells.Item(45, 6).Value = 371
$ws.Cells.Item(48, 6).Value = 4122
$ws.Cells.Item(49, 6).Value = 151

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 4173
$ws.Cells.Item(4, 6).Value = 70
$ws.Cells.Item(5, 6).Value = 104
$ws.Cells.Item(7, 6).Value = 47
$ws.Cells.Item(12, 6).Value = 1055
$ws.Cells.Item(18, 6).Value = 52

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 766
$ws.Cells.Item(3, 6).Value = 489
$ws.Cells.Item(4, 6).Value = 107
$ws.Cells.Item(5, 6).Value = 18

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 489
$ws.Cells.Item(3, 6).Value = 155
$ws.Cells.Item(4, 6).Value = 929
$ws.Cells.Item(5, 6).Value = 1090
$ws.Cells.Item(6, 6).Value = 1554
$ws.Cells.Item(7, 6).Value = 339
$ws.Cells.Item(8, 6).Value = 697
$ws.Cells.Item(9, 6).Value = 12532
$ws.Cells.Item(10, 6).Value = 2210
$ws.Cells.Item(12, 6).Value = 275
$ws.Cells.Item(13, 6).Value = 1252
$ws.Cells.Item(14, 6).Value = 238
$ws.Cells.Item(16, 6).Value = 801
$ws.Cells.Item(18, 6).Value = 319
$ws.Cells.Item(20, 6).Value = 788
$ws.Cells.Item(21, 6).Value = 4475
$ws.Cells.Item(22, 6).Value = 4475
$ws.Cells.Item(23, 6).Value = 1158
$ws.Cells.Item(24, 6).Value = 18
$ws.Cells.Item(25, 6).Value = 104
$ws.Cells.Item(26, 6).Value = 13
$ws.Cells.Item(27, 6).Value = 47
$ws.Cells.Item(31, 6).Value = 1104
$ws.Cells.Item(33, 6).Value = 120
$ws.Cells.Item(35, 6).Value = 283
$ws.Cells.Item(39, 6).Value = 4614
$ws.Cells.Item(41, 6).Value = 179
$ws.Cells.Item(46, 6).Value = 4122
$ws.Cells.Item(47, 6).Value = 52
